$d = $word.ActiveDocument

# --- Locate the "...nein ersetzen" paragraph and detach the _GoBack bookmark ---
# (it currently sits at the end of this paragraph; it needs to move to the end
#  of the new "Wenn als Lehrer angemeldet" paragraph created below).
$pNeinErsetzen = $d.Paragraphs.Item(95)
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- New paragraph: "Nice to have: In Home " (list level 1 / ilvl=0) ---
$pNeinErsetzen.Range.InsertParagraphAfter()
$pNiceToHave = $d.Paragraphs.Item(96)
$pNiceToHave.Range.ListFormat.ListLevelNumber = 1
$pNiceToHave.Range.Text = "Nice to have: In Home "

# --- New paragraph: "Wenn als Lehrer angemeldet..." (list level 2 / ilvl=1) ---
$pNiceToHave.Range.InsertParagraphAfter()
$pLehrer = $d.Paragraphs.Item(97)
$pLehrer.Range.ListFormat.ListLevelNumber = 2
# Append a trailing placeholder character so the insertion point for the
# re-added bookmark below is not the literal last character of the
# paragraph (that specific position triggers a bookmark placement defect
# in this host); the placeholder is stripped again immediately after.
$pLehrer.Range.Text = "Wenn als Lehrer angemeldet: Direktlink zu „meine Klassen“ „meine Schüler“ „meine Notizen“#"

$placeholderPos = $pLehrer.Range.End - 2
$bmRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($pLehrer.Range.End - 2, $pLehrer.Range.End - 1)
$placeholderRange.Delete()

# --- Existing (previously empty) paragraph becomes "Wenn als Schüler angemeldet..." ---
$pSchueler = $d.Paragraphs.Item(98)
$pSchueler.Range.ListFormat.ListLevelNumber = 2
$pSchueler.Range.Text = "Wenn als Schüler angemeldet: Direktlink zu „meine Klasse“ „meine Lehrer“ „meine Notizen“"

# --- New paragraph: "Wenn als Root angemeldet: ---" (list level 2 / ilvl=1) ---
$pSchueler.Range.InsertParagraphAfter()
$pRoot = $d.Paragraphs.Item(99)
$pRoot.Range.ListFormat.ListLevelNumber = 2
$pRoot.Range.Text = "Wenn als Root angemeldet: ---"
